$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.118.30"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "3.588.18"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'208.12"
$ws.Range("E5").Value = "  +7.53%  "
$ws.Range("D6").Value = "'569.23"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.682"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'63.27"
$ws.Range("E10").Value = "  +13.07%  "
$ws.Range("D11").Value = "'0.147"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "'10.36"
$ws.Range("E13").Value = "  +4.95%  "
$ws.Range("D14").Value = "4.161.63"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "3.588.29"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'19.19"
$ws.Range("E16").Value = "  +4.05%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.126"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "67.936.21"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "'12.20"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'402.88"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").Value = "'12.34"
$ws.Range("E23").Value = "  +7.95%  "
$ws.Range("D24").Value = "'84.62"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'12.51"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'3.86"
$ws.Range("E27").Value = "  +5.69%  "
$ws.Range("D28").Value = "'9.28"
$ws.Range("E28").Value = "  +3.38%  "
$ws.Range("D29").Value = "'7.62"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'31.53"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "'688.33"
$ws.Range("E31").Value = "  +9.04%  "
$ws.Range("D32").Value = "'12.15"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "'63.34"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").Value = "'41.38"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "  +7.58%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.17"
$ws.Range("E39").Value = "  +20.76%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0747"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.159.71"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.133"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'2.77"
$ws.Range("E45").Value = "  +8.72%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0412"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").Value = "'3.13"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").Value = "'8.69"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").Value = "'138.94"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("E51").Value = "  -1.80%  "
